# Oct 2019 to Nov 2019 - Cord Cutting Wizard - Change Log
# Reword the nine "Comment" category labels used in column E (the wording
# was tightened/recapitalized; the category each row belongs to is
# unchanged) and move the active selection to E5:E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oct 2019 to Nov 2019")

$commentRange = $ws.Range("E2:E59")

# Old label -> new label, applied as whole-cell replacements so partial
# text overlaps between old/new phrases can't cross-contaminate rows.
$commentRange.Replace("New Alias added in Nov 2019", "Alias Changed for Network", 2)
$commentRange.Replace("Network removed from base Service in Nov 2019", "Network Removed from Base Service", 2)
$commentRange.Replace("Network moved from base Service to Add-On Package in Nov 2019", "Network Moved from Base Service to Add-On Package in Nov 2019", 2)
$commentRange.Replace("Network added to Base Service in Nov 2019", "Network Added to Base Service in Nov 2019", 2)
$commentRange.Replace("Network moved from Add-On Package to Base Service in Nov 2019", "Network Moved from Add-On Package to Base Service in Nov 2019", 2)
$commentRange.Replace("New Network added to database in Nov 2019", "New Network Added to Database in Nov 2019", 2)
$commentRange.Replace("Network removed from database in Nov 2019", "Network Removed from Database", 2)
$commentRange.Replace("Add-On Package renamed in Nov 2019", "Name of Add-On Package Changed", 2)
$commentRange.Replace("Network added to Add-On Package in Nov 2019", "Network Added to Add-On Package", 2)

# Move the selection as recorded in the saved view state.
$ws.Range("E5:E6").Select()
